$wb = $excel.ActiveWorkbook

# Rename existing "eia_area" sheet to "eia_location_id"
$wsLoc = $wb.Worksheets.Item("eia_area")
$wsLoc.Name = "eia_location_id"

# Add a new sheet "eia_area" right after "eia_location_id"
$wsArea = $wb.Worksheets.Add($null, $wsLoc)
$wsArea.Name = "eia_area"

# Populate header row in C, D, B, A order (matches new shared-string order)
$wsArea.Range("C1").Value = "sqft"
$wsArea.Range("D1").Value = "Impervious_sqft"
$wsArea.Range("B1").Value = "acres"
$wsArea.Range("A1").Value = "location_id"

# Populate BMP data rows (location_id / acres / sqft / Impervious_sqft)
$wsArea.Range("A2").Value = "TPO01"
$wsArea.Range("B2").Value = 71.444418330000005
$wsArea.Range("C2").Value = 5906012.3820000002
$wsArea.Range("D2").Value = 2675408.4810000001
$wsArea.Range("A3").Value = "TDO01"
$wsArea.Range("B3").Value = 136.07334739999999
$wsArea.Range("C3").Value = 5927331.3049999997
$wsArea.Range("D3").Value = 1625421.122
$wsArea.Range("A4").Value = "TKV01"
$wsArea.Range("B4").Value = 171.85245789999999
$wsArea.Range("C4").Value = 7485863.1229999997
$wsArea.Range("D4").Value = 2819719.7179999999
$wsArea.Range("A5").Value = "TMH01"
$wsArea.Range("B5").Value = 174.32630649999999
$wsArea.Range("C5").Value = 7593623.5350000001
$wsArea.Range("D5").Value = 2576687.5819999999
$wsArea.Range("A6").Value = "TNS01"
$wsArea.Range("B6").Value = 217.7788955
$wsArea.Range("C6").Value = 9486410.7420000006
$wsArea.Range("D6").Value = 3660210.128
$wsArea.Range("A7").Value = "TBK01"
$wsArea.Range("B7").Value = 231.58095510000001
$wsArea.Range("C7").Value = 10087626.050000001
$wsArea.Range("D7").Value = 2137233.2710000002
$wsArea.Range("A8").Value = "TFD01"
$wsArea.Range("B8").Value = 234.09865740000001
$wsArea.Range("C8").Value = 10197296.73
$wsArea.Range("D8").Value = 3478720.3960000002
$wsArea.Range("A9").Value = "TTX27"
$wsArea.Range("B9").Value = 249.2998709
$wsArea.Range("C9").Value = 10859458.939999999
$wsArea.Range("D9").Value = 4057289.39
$wsArea.Range("A10").Value = "TFS01"
$wsArea.Range("B10").Value = 276.8609783
$wsArea.Range("C10").Value = 12060015.970000001
$wsArea.Range("D10").Value = 5035429.9689999996
$wsArea.Range("A11").Value = "TDA01"
$wsArea.Range("B11").Value = 270.26444609999999
$wsArea.Range("C11").Value = 12104505.98
$wsArea.Range("D11").Value = 3150624.395
$wsArea.Range("A12").Value = "TPB01"
$wsArea.Range("B12").Value = 280.45176459999999
$wsArea.Range("C12").Value = 12216430
$wsArea.Range("D12").Value = 3607941.4730000002
$wsArea.Range("A13").Value = "TFC01"
$wsArea.Range("B13").Value = 291.45618960000002
$wsArea.Range("C13").Value = 12695780.83
$wsArea.Range("D13").Value = 4527052.1179999998
$wsArea.Range("A14").Value = "TDU01"
$wsArea.Range("B14").Value = 409.25697700000001
$wsArea.Range("C14").Value = 17827162.609999999
$wsArea.Range("D14").Value = 1609838.871
$wsArea.Range("A15").Value = "TSO01"
$wsArea.Range("B15").Value = 516.37507649999998
$wsArea.Range("C15").Value = 22493208.359999999
$wsArea.Range("D15").Value = 10075179
$wsArea.Range("A16").Value = "TNA01"
$wsArea.Range("B16").Value = 472.15890439999998
$wsArea.Range("C16").Value = 28262408.210000001
$wsArea.Range("D16").Value = 10428891.060000001
$wsArea.Range("A17").Value = "TLU01"
$wsArea.Range("B17").Value = 643.53066860000001
$wsArea.Range("C17").Value = 28811482.870000001
$wsArea.Range("D17").Value = 14779080.73
$wsArea.Range("A18").Value = "TPI01"
$wsArea.Range("B18").Value = 449.3504049
$wsArea.Range("C18").Value = 28934650.41
$wsArea.Range("D18").Value = 7336519.3689999999
$wsArea.Range("A19").Value = "TFE01"
$wsArea.Range("B19").Value = 220.8060729
$wsArea.Range("C19").Value = 41747880.719999999
$wsArea.Range("D19").Value = 18710988.010000002
$wsArea.Range("A20").Value = "THR01"
$wsArea.Range("B20").Value = 1103.8468869999999
$wsArea.Range("C20").Value = 48083378.079999998
$wsArea.Range("D20").Value = 20152081.27
$wsArea.Range("A21").Value = "TFB01"
$wsArea.Range("B21").Value = 1134.548792
$wsArea.Range("C21").Value = 49420747.68
$wsArea.Range("D21").Value = 18262208.25
$wsArea.Range("A22").Value = "TBR01"
$wsArea.Range("B22").Value = 1148.451286
$wsArea.Range("C22").Value = 51304948.560000002
$wsArea.Range("D22").Value = 17934581.789999999
$wsArea.Range("A23").Value = "RCR01"
$wsArea.Range("B23").Value = 1603.223352
$wsArea.Range("C23").Value = 69836129.849999994
$wsArea.Range("D23").Value = 50051478.759999998
$wsArea.Range("A24").Value = "TWB06"
$wsArea.Range("B24").Value = 848.33045019999997
$wsArea.Range("C24").Value = 90723670.299999997
$wsArea.Range("D24").Value = 32039217.879999999
$wsArea.Range("A25").Value = "TWB05"
$wsArea.Range("B25").Value = 410.63448010000002
$wsArea.Range("C25").Value = 108610836.7
$wsArea.Range("D25").Value = 38050077.590000004
$wsArea.Range("A26").Value = "TPY01"
$wsArea.Range("B26").Value = 158.9618614
$wsArea.Range("C26").Value = 109878251
$wsArea.Range("D26").Value = 57375195.619999997
$wsArea.Range("A27").Value = "TOR01"
$wsArea.Range("B27").Value = 2171.0271550000002
$wsArea.Range("C27").Value = 414641447.69999999
$wsArea.Range("D27").Value = 145300487.90000001
$wsArea.Range("A28").Value = "RCR05"
$wsArea.Range("B28").Value = 1335.459206
$wsArea.Range("C28").Value = 797091117.29999995
$wsArea.Range("D28").Value = 266505068.69999999
$wsArea.Range("A29").Value = "RCR09"
$wsArea.Range("B29").Value = 925.03015310000001
$wsArea.Range("C29").Value = 1033762101
$wsArea.Range("D29").Value = 393938764.89999998

# Approximate the target column widths (engine applies its own px rounding)
$wsArea.Columns.Item(1).ColumnWidth = 6.6328125
$wsArea.Range("B1:C1").EntireColumn.ColumnWidth = 11.81640625
$wsArea.Columns.Item(4).ColumnWidth = 19.26953125

# Restore the selection on the now-inactive eia_location_id sheet
$wsLoc.Range("E29").Select()

# eia_area becomes the active/selected sheet with this cell selected
$wsArea.Activate()
$wsArea.Range("J25").Select()
